$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.0002567018593964017
$ws.Range("E2").Value = 0.0002567018593964017

# Row 3
$ws.Range("D3").Value = 0.9186475293880044
$ws.Range("E3").Value = 0.9186475293880044

# Row 4
$ws.Range("D4").Value = 0.0005943758585822261
$ws.Range("E4").Value = 0.0005943758585822261

# Row 5
$ws.Range("D5").Value = 0.0006934940963239417
$ws.Range("E5").Value = 0.0006934940963239417

# Row 6
$ws.Range("D6").Value = 0.1143427434842015
$ws.Range("E6").Value = 0.1143427434842015

# Row 7
$ws.Range("D7").Value = 0.8866943030131962
$ws.Range("E7").Value = 0.1133056969868038

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = [double]"6.423438187741829E-05"
$ws.Range("E8").Value = 0.9999357656181226

# Row 9
$ws.Range("D9").Value = 0.9876016543609304
$ws.Range("E9").Value = 0.01239834563906961

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 0.001106210598497721
$ws.Range("E10").Value = 0.9988937894015023

# Row 11
$ws.Range("D11").Value = 0.9999999999558611
$ws.Range("E11").Value = [double]"4.413891474541742E-11"
$ws.Range("F11").Value = 1.922445058822632
$ws.Range("G11").Value = 0.7

# Row 12
$ws.Range("D12").Value = 0.001070144132705447
$ws.Range("E12").Value = 0.001070144132705447

# Row 13
$ws.Range("D13").Value = 0.9996322199442315
$ws.Range("E13").Value = 0.9996322199442315

# Row 14
$ws.Range("D14").Value = 0.0007198397294128899
$ws.Range("E14").Value = 0.0007198397294128899

# Row 15
$ws.Range("D15").Value = 0.001677314130025276
$ws.Range("E15").Value = 0.001677314130025276

# Row 16
$ws.Range("D16").Value = 0.03077223024567848
$ws.Range("E16").Value = 0.03077223024567848

# Row 17
$ws.Range("D17").Value = 0.9671782497501139
$ws.Range("E17").Value = 0.03282175024988609

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = 0.0005039164600994533
$ws.Range("E18").Value = 0.9994960835399005

# Row 19
$ws.Range("D19").Value = 0.9537390256791032
$ws.Range("E19").Value = 0.04626097432089682

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = 0.003238113777760553
$ws.Range("E20").Value = 0.9967618862222395

# Row 21
$ws.Range("D21").Value = 0.9999999999999993
$ws.Range("E21").Value = [double]"6.661338147750939E-16"
$ws.Range("F21").Value = 2.13493537902832
$ws.Range("G21").Value = 0.7
